# live_trading_results.xlsx update
# Trade #51 closed at 2026-02-18 00:20:39 - unknown UNKNOWN +0.000%
# and a new trade (#109) logged for the momentum strategy.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet: refresh headline stats
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.47   # Current Capital
$wsSummary.Range("B4").Value = 0.58      # Total P&L $
$wsSummary.Range("B6").Value = 79        # Total Trades
$wsSummary.Range("B8").Value = 32        # Losing Trades
$wsSummary.Range("B9").Value = 49.37     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet: "momentum" row (row 11) stats
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Cells.Item(11, 3).Value = 99.63   # Capital
$wsStatus.Cells.Item(11, 4).Value = 11      # Trades
$wsStatus.Cells.Item(11, 5).Value = -0.37   # P&L $
$wsStatus.Cells.Item(11, 6).Value = -0.37   # P&L %
$wsStatus.Cells.Item(11, 7).Value = 9.09    # Win Rate %

# ---------------------------------------------------------------
# All Trades sheet: close out trade #80 (row 81) and append #109
# Columns: A Trade#, B Date, C Time, D Strategy, E Side,
#          F Entry Price, G Exit Price, H Status, I P&L%, J P&L$,
#          K Capital After, L Exit Reason, M Duration (min),
#          N Entry Slippage (bps), O Exit Slippage (bps),
#          P Confidence, Q Entry Reason
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Cells.Item(81, 7).Value = 0.65          # Exit Price
$wsAll.Cells.Item(81, 8).Value = "CLOSED"      # Status
$wsAll.Cells.Item(81, 9).Value = -2.9851       # P&L %
$wsAll.Cells.Item(81, 10).Value = -0.02        # P&L $
$wsAll.Cells.Item(81, 11).Value = 99.63        # Capital After
$wsAll.Cells.Item(81, 12).Value = "early_exit" # Exit Reason
$wsAll.Cells.Item(81, 13).Value = 0.12         # Duration (min)

$wsAll.Cells.Item(110, 1).Value = 109
$wsAll.Cells.Item(110, 2).NumberFormat = "@"
$wsAll.Cells.Item(110, 2).Value = "2026-02-18"
$wsAll.Cells.Item(110, 2).Style = "Normal"
$wsAll.Cells.Item(110, 3).Value = "00:20:33"
$wsAll.Cells.Item(110, 4).Value = "momentum"
$wsAll.Cells.Item(110, 5).Value = "DOWN"
$wsAll.Cells.Item(110, 6).Value = 0.67
$wsAll.Cells.Item(110, 11).Value = 99.64873713109129
$wsAll.Cells.Item(110, 16).Value = 0.9
$wsAll.Cells.Item(110, 17).Value = "Downward momentum: -3.810% over 10 samples"
$wsAll.Cells.Item(110, 8).Value = "OPEN"
$wsAll.Cells.Item(110, 9).Value = 0
$wsAll.Cells.Item(110, 10).Value = 0
$wsAll.Cells.Item(110, 13).Value = 0
$wsAll.Cells.Item(110, 14).Value = 0
$wsAll.Cells.Item(110, 15).Value = 0

# ---------------------------------------------------------------
# momentum strategy sheet: mirror the same trade #80 close-out
# (row 13) and append trade #109 (row 28).
# Columns: A Trade#, B Date, C Time, D Strategy, E Side,
#          F Entry Price, G Exit Price, H Status, I P&L%, J P&L$,
#          K Capital After, L Entry Slippage (bps),
#          M Exit Slippage (bps), N Confidence, O Entry Reason,
#          P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------
$wsMom = $wb.Worksheets.Item("momentum")
$wsMom.Cells.Item(13, 7).Value = 0.65            # Exit Price
$wsMom.Cells.Item(13, 8).Value = "CLOSED"        # Status
$wsMom.Cells.Item(13, 9).Value = -2.9851         # P&L %
$wsMom.Cells.Item(13, 10).Value = -0.02          # P&L $
$wsMom.Cells.Item(13, 11).Value = 99.63          # Capital After
$wsMom.Cells.Item(13, 16).Value = "early_exit"   # Exit Reason
$wsMom.Cells.Item(13, 17).Value = 0.12           # Duration (min)

$wsMom.Cells.Item(28, 1).Value = 109
$wsMom.Cells.Item(28, 2).NumberFormat = "@"
$wsMom.Cells.Item(28, 2).Value = "2026-02-18"
$wsMom.Cells.Item(28, 2).Style = "Normal"
$wsMom.Cells.Item(28, 3).Value = "00:20:33"
$wsMom.Cells.Item(28, 4).Value = "momentum"
$wsMom.Cells.Item(28, 5).Value = "DOWN"
$wsMom.Cells.Item(28, 6).Value = 0.67
$wsMom.Cells.Item(28, 8).Value = "OPEN"
$wsMom.Cells.Item(28, 9).Value = 0
$wsMom.Cells.Item(28, 10).Value = 0
$wsMom.Cells.Item(28, 11).Value = 99.64873713109129
$wsMom.Cells.Item(28, 12).Value = 0
$wsMom.Cells.Item(28, 13).Value = 0
$wsMom.Cells.Item(28, 14).Value = 0.9
$wsMom.Cells.Item(28, 15).Value = "Downward momentum: -3.810% over 10 samples"
$wsMom.Cells.Item(28, 17).Value = 0
